# "Register by email, small ui fixes"
# Adds an "email" column to the survey users table and records a new
# registrant (Natalija), while refreshing Aleksandar's timestamp/code too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "email" column between "last_name" and "timestamp".
$ws.Columns("C").Insert()
$ws.Range("C1").Value = "email"

# Row 2: new registrant Natalija Gajic.
$ws.Range("A2").Value = "Natalija"
$ws.Range("B2").Value = "Gajic"
$ws.Range("C2").Value = "2023_07_11_21_27_31"
$ws.Range("D2").Value = "nat.gaj98@gmail.com"
$ws.Range("E2").Value = "1lbCvs1hADW3mD5c66mzpqb7ziKmDGtr5QoGX7pISYs"
$ws.Range("F2").Value = $false

# Row 3: Aleksandar Gajic, now with refreshed timestamp/code and his email.
$ws.Range("A3").Value = "Aleksandar"
$ws.Range("B3").Value = "Gajic"
$ws.Range("C3").Value = "2023_07_11_21_31_34"
$ws.Range("D3").Value = "gajic7080@gmail.com"
$ws.Range("E3").Value = "Qmi9nT6p5G-tm03YAAsHhCaf_5XHOwXocy9IiDPYSKc"
$ws.Range("F3").Value = $false
